$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.48602283000946
$ws.Range("B1").Value = 3.922909021377563
$ws.Range("C1").Value = 3.676296234130859
$ws.Range("D1").Value = 1.544712066650391
$ws.Range("E1").Value = 0.990425169467926
